$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 24 (CRP_nr): columns C-H correspond to years 1987,1992,1997,2002,2007,2012
$ws.Range("C24").Value = 92.0228271484375
$ws.Range("D24").Value = 76.53497314453125
$ws.Range("E24").Value = 67.908920288085938
$ws.Range("F24").Value = 60.107421875
$ws.Range("G24").Value = 55.856372833251953
$ws.Range("H24").Value = 55.617328643798828

# Row 25 (pasture_nr): columns E-H correspond to years 1997,2002,2007,2012
$ws.Range("E25").Value = 18.795343399047852
$ws.Range("F25").Value = 17.547082901000977
$ws.Range("G25").Value = 18.402990341186523
$ws.Range("H25").Value = 16.522886276245117
